$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 is no longer the last row of data, so it should take on the
# same date format as the rest of the body rows (A2:A37) rather than
# the special "last row" format.
$ws.Range("A38").NumberFormat = $ws.Range("A37").NumberFormat

# Append the new day's data as row 39.
$ws.Range("A39").Value = 45623
$ws.Range("B39").Value = 102
$ws.Range("C39").Value = 83
$ws.Range("D39").Value = 93

# Row 39 is now the last row, so it gets the "last row" date format
# that row 38 used to have.
$ws.Range("A39").NumberFormat = "YYYY-MM-DD"
